# Update cryptos list with latest prices / volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.506.33'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -2.67%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.355.90'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -4.55%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '555.80'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -4.87%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '176.04'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.93%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.617'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -3.20%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.348.33'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -4.42%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.628'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.08%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.163'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.10%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.98'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.22%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000273'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -3.12%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.05'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.05%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.898.17'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -4.41%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.69%  '
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.387.85'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -3.72%  '
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.118'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -2.77%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.82'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -2.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '64.500.20'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.69%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.980'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -3.61%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '437.42'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +5.18%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +12.17%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -5.00%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.58'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.91%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.47'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.30%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.99%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.64%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.76'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -4.88%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '29.71'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.62%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.63'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.60%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.46'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '578.41'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.27%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -3.34%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '58.43'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -3.81%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.07%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -8.22%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.51'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -4.09%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '35.72'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -3.42%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0₃0754'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -6.16%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -5.12%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.101.76'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -4.28%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.11%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.81'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -5.61%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -3.66%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0410'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -3.06%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -4.30%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.59%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -3.68%  '
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '136.96'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.29%  '
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'THORChain'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.29'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -4.32%  '
